$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Main": correct the NAV for 2022-07-25 (row 11) and extend the
# series with a new row (12) for the next trading day.
# ---------------------------------------------------------------------------
$main = $wb.Worksheets.Item("Main")

# Update C11 with corrected NAV value
$main.Range("C11").Value = 52314.02

# Fill in D11:G11 to match the established row formulas
$main.Range("D11").Formula = "=C11-C10"
$main.Range("F11").NumberFormat = $main.Range("F10").NumberFormat
$main.Range("F11").Formula = "=+C11/C10-1"

# Add new row 12 for the next date
$main.Range("B12").Formula = "=+B11+1"

[void]$main.Range("C13").Select()

# ---------------------------------------------------------------------------
# Sheet "Trades": append newly executed trades from 2022-07-26.
# ---------------------------------------------------------------------------
$trades = $wb.Worksheets.Item("Trades")
$trades.Activate()
$excel.ActiveWindow.ScrollRow = 183
$excel.ActiveWindow.ScrollColumn = 3

function Set-Row {
    param($r, $b, $c, $d, $e, $f, $g, $h, $i, $j, $k, $l)
    if ($b -ne $null) { $trades.Range("B$r").Value = $b }
    if ($c -ne $null) { $trades.Range("C$r").Value = $c }
    if ($d -ne $null) { $trades.Range("D$r").Value = $d }
    if ($e -ne $null) { $trades.Range("E$r").Value = $e }
    if ($f -ne $null) { $trades.Range("F$r").Value = $f }
    if ($g -ne $null) { $trades.Range("G$r").Value = $g }
    if ($h -ne $null) { $trades.Range("H$r").Value = $h }
    if ($i -ne $null) { $trades.Range("I$r").Value = $i }
    if ($j -ne $null) { $trades.Range("J$r").Value = $j }
    if ($k -ne $null) { $trades.Range("K$r").Value = $k }
    if ($l -ne $null) {
        # "Total <symbol>" summary rows carry the same numeric style (#,##0.00)
        # on column L as the rest of the row, even though L just holds a
        # placeholder space string there.
        if ($l -eq " ") { $trades.Range("L$r").NumberFormat = "#,##0.00" }
        $trades.Range("L$r").Value = $l
    }
}

Set-Row 190 "AMC"       "2022-07-26, 09:30:58" 100  14.585 14.03  -1458.5   -1     1459.5    0       -55.5   "O"
Set-Row 191 "Total AMC" $null                  100  " "    $null  -1458.5   -1     1459.5    0       -55.5   $null

Set-Row 192 "AVXL"       "2022-07-26, 13:50:44" 300  11.048 10.79  -3314.4   -1.5   3445.38   129.48  -77.4   "C"
Set-Row 193 "Total AVXL" $null                  300  " "    $null  -3314.4   -1.5   3445.38   129.48  -77.4   " "

Set-Row 194 "SHOP" "2022-07-26, 09:36:54" 100  30.92       31.55  -3092     -1     3093      0       63     "O"
Set-Row 195 "SHOP" "2022-07-26, 09:38:07" 100  31.19       31.55  -3119     -1     3120      0       36     "O"
Set-Row 196 "SHOP" "2022-07-26, 13:22:20" -50  31.252      31.55  1562.6    -1.04  -1546.5   15.06   -14.9  "C;P"
Set-Row 197 "SHOP" "2022-07-26, 13:49:34" -150 31.215      31.55  4682.25   -1.13  -4666.5   14.62   -50.25 "C"
Set-Row 198 "Total SHOP" $null             0    " "        $null  33.85     -4.17  0         29.68   33.85  " "

Set-Row 199 "SIGA" "2022-07-26, 09:55:17" -100 15.9        17.43  1590      -1.05  -1588.95  0       -153   "O"
Set-Row 200 "SIGA" "2022-07-26, 10:02:10" -100 16.31       17.43  1631      -1.05  -1629.95  0       -112   "O"
Set-Row 201 "SIGA" "2022-07-26, 10:05:43" -100 16.58       17.43  1658      -1.05  -1656.95  0       -85    "O;P"
Set-Row 202 "SIGA" "2022-07-26, 13:50:51" -100 17.5        17.43  1750      -1.05  -1748.95  0       7      "O;P"
Set-Row 203 "Total SIGA" $null             -400 " "        $null  6629      -4.2   -6624.8   0       -343   " "

Set-Row 204 "SPY" "2022-07-26, 09:51:41" -50  393.46      390.89 19673     -1.46  -19671.54 198.61  128.5  "O"
Set-Row 205 "SPY" "2022-07-26, 10:58:28" 10   392.65      390.89 -3926.5   -1     3934.31   6.81    -17.6  "C"
Set-Row 206 "SPY" "2022-07-26, 11:04:01" 10   392.6195    390.89 -3926.2   -1     3934.31   7.11    -17.3  "C"
Set-Row 207 "SPY" "2022-07-26, 12:14:00" 10   390.5062    390.89 -3905.06  -1     3886.68   -19.38  3.84   "C"
Set-Row 208 "SPY" "2022-07-26, 12:45:32" 20   390.966     390.89 -7819.32  -1     7717.64   -102.68 -1.52  "C"
Set-Row 209 "Total SPY" $null             0    " "        $null  95.92     -5.46  -198.61   90.47   95.92  " "

Set-Row 210 "TWTR" "2022-07-26, 09:31:37" 100  39.08       39.34  -3908     -1     3902.33   2.3     26    "C"
Set-Row 211 "TWTR" "2022-07-26, 09:38:49" 100  38.94       39.34  -3894     -1     3908.18   13.18   40    "C"
Set-Row 212 "TWTR" "2022-07-26, 12:45:40" -100 39.34       39.34  3934      -1.1   -3932.9   0       0     "O"
Set-Row 213 "TWTR" "2022-07-26, 13:49:22" 200  39.258      39.34  -7851.6   -1     7871.07   18.47   16.4  "C"
Set-Row 214 "Total TWTR" $null             300  " "        $null  -11719.6  -4.1   11748.69  33.95   82.4  " "

Set-Row 215 "Total" $null $null $null $null -9733.73 -20.43 9830.17 283.58 -263.73 " "
$trades.Range("H215").Style = "Normal"

[void]$trades.Range("A216").Select()
